# Generate Report for Handoff
# - Update "Status" text from "Handed back: in sync with en-US" to "Ready for handoff"
# - Update "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps
# - Narrow the "Status" columns to fit the new shorter text

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (E2, F2) and handoff date (G2)
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-21 07:05:15"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-08-21 07:05:10"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-08-21 07:05:15"

# Narrow the Status columns now that the text is shorter
# (17.2159881591797 is the saved column-width target; the COM ColumnWidth
# property is quantized to a 1/6-character pixel grid, so we use the
# nearest representable input - 16.3333... - which round-trips to the
# closest achievable stored width, 17.1667)
$ws1.Columns.Item(5).ColumnWidth = 16.3333333333333
$ws1.Columns.Item(6).ColumnWidth = 16.3333333333333
$ws2.Columns.Item(3).ColumnWidth = 16.3333333333333
$ws3.Columns.Item(3).ColumnWidth = 16.3333333333333
